$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to Text
# so Excel smart type-detection does not silently convert "1.00" -> 1, etc.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '62.789.77'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.439.29'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '567.01'
$ws.Range('E5').Value = '  +0.78%  '
Set-TextValue 'D6' '145.70'
$ws.Range('E6').Value = '  +2.14%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('E10').Value = '  +0.56%  '
Set-TextValue 'D11' '5.25'
$ws.Range('E11').Value = '  -0.99%  '
Set-TextValue 'D12' '0.353'
$ws.Range('E12').Value = '  +0.70%  '
Set-TextValue 'D13' '0.0000186'
$ws.Range('E13').Value = '  +7.43%  '
Set-TextValue 'D14' '26.87'
$ws.Range('E14').Value = '  +4.91%  '
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '62.634.86'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Value = '2.443.47'
$ws.Range('E17').Value = '  +1.21%  '
Set-TextValue 'D18' '11.26'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('E19').Value = '  +1.57%  '
Set-TextValue 'D20' '324.07'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('E22').Value = '  -0.06%  '
Set-TextValue 'D23' '67.22'
$ws.Range('E23').Value = '  +1.70%  '
Set-TextValue 'D24' '1.79'
$ws.Range('E24').Value = '  +3.75%  '
Set-TextValue 'D25' '8.75'
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('E26').Value = '  +8.69%  '
Set-TextValue 'D27' '572.38'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').Value = '2.559.05'
$ws.Range('E28').Value = '  +1.23%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D30' '8.44'
$ws.Range('E30').Value = '  +2.99%  '
Set-TextValue 'D31' '1.46'
$ws.Range('E31').Value = '  +3.08%  '
Set-TextValue 'D32' '0.148'
$ws.Range('E32').Value = '  -0.25%  '
Set-TextValue 'D33' '1.87'
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  +1.37%  '
Set-TextValue 'D35' '4.90'
$ws.Range('E35').Value = '  +4.37%  '
$ws.Range('E36').Value = '  -0.16%  '
Set-TextValue 'D37' '0.384'
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('E38').Value = '  -0.57%  '
Set-TextValue 'D39' '18.81'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D40' '1.84'
$ws.Range('E40').Value = '  +2.63%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D41' '148.63'
$ws.Range('E41').Value = '  -1.86%  '
$ws.Range('E42').Value = '  +0.16%  '
Set-TextValue 'D43' '2.45'
$ws.Range('E43').Value = '  +7.15%  '
Set-TextValue 'D44' '149.32'
$ws.Range('E44').Value = '  +0.96%  '
Set-TextValue 'D45' '3.68'
$ws.Range('E45').Value = '  +1.49%  '
Set-TextValue 'D46' '0.0538'
$ws.Range('E46').Value = '  +1.09%  '
Set-TextValue 'D47' '20.67'
$ws.Range('E47').Value = '  +3.77%  '
Set-TextValue 'D48' '0.601'
$ws.Range('E48').Value = '  +1.39%  '
Set-TextValue 'D49' '0.0232'
$ws.Range('E49').Value = '  +2.82%  '
Set-TextValue 'D50' '0.0928'
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D51' '11.61'
$ws.Range('E51').Value = '  +0.66%  '
